$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 3231.6667
$ws.Range("I12").Value = 1558.8889
$ws.Range("K12").Value = 1558.8889
$ws.Range("M12").Value = -1388.8889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 878.4167
$ws.Range("I33").Value = 871.86365
$ws.Range("J33").Value = 950.5
$ws.Range("K33").Value = 871.86365
$ws.Range("L33").Value = 950.5
$ws.Range("M33").Value = -642.86365
$ws.Range("N33").Value = -1408.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4501
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7500.143
$ws.Range("I51").Value = 6832.6665
$ws.Range("J51").Value = 8000.75
$ws.Range("K51").Value = 6832.6665
$ws.Range("L51").Value = 8000.75
$ws.Range("M51").Value = -6348.6665
$ws.Range("N51").Value = -8968.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 31252038
$ws.Range("J62").Value = 2350
$ws.Range("L62").Value = 2350
$ws.Range("N62").Value = -3598

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 9166.666999999999
$ws.Range("I64").Value = 6000
$ws.Range("J64").Value = 15500
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 15500
$ws.Range("M64").Value = -5752
$ws.Range("N64").Value = -15996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 31252038
$ws.Range("J65").Value = 2350
$ws.Range("L65").Value = 11750
$ws.Range("N65").Value = -17990

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 9166.666999999999
$ws.Range("I67").Value = 6000
$ws.Range("J67").Value = 15500
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 15500
$ws.Range("M67").Value = -5142
$ws.Range("N67").Value = -17216

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5265887.5
$ws.Range("I86").Value = 3599.4
$ws.Range("K86").Value = 3599.4
$ws.Range("M86").Value = -2476.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 5265887.5
$ws.Range("I89").Value = 3599.4
$ws.Range("K89").Value = 17997
$ws.Range("M89").Value = -12381

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 450.1875
$ws.Range("J96").Value = 649.4
$ws.Range("L96").Value = 1948.2
$ws.Range("N96").Value = -4694.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1750.3
$ws.Range("I98").Value = 1437.8379
$ws.Range("K98").Value = 1437.8379
$ws.Range("M98").Value = 60.16210000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1982.7778
$ws.Range("I106").Value = 3954.6667
$ws.Range("K106").Value = 3954.6667
$ws.Range("M106").Value = -3323.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4947.0386
$ws.Range("J112").Value = 5560.136
$ws.Range("L112").Value = 16680.408
$ws.Range("N112").Value = -18896.408

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1750.3
$ws.Range("I122").Value = 1437.8379
$ws.Range("K122").Value = 4313.5137
$ws.Range("M122").Value = -1863.5137

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6321.9873
$ws.Range("I138").Value = 3422.5
$ws.Range("J138").Value = 7771.731
$ws.Range("K138").Value = 10267.5
$ws.Range("L138").Value = 23315.193
$ws.Range("M138").Value = -5127.5
$ws.Range("N138").Value = -33595.193

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 85747.46000000001
$ws.Range("I2").Value = 11922.889
$ws.Range("K2").Value = 11922.889
$ws.Range("M2").Value = -11809.889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2661.9214
$ws.Range("I32").Value = 2159.8965
$ws.Range("K32").Value = 2159.8965
$ws.Range("M32").Value = -1872.8965

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1688.3158
$ws.Range("I74").Value = 1726.5555
$ws.Range("K74").Value = 1726.5555
$ws.Range("M74").Value = -852.5554999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1688.3158
$ws.Range("I77").Value = 1726.5555
$ws.Range("K77").Value = 8632.7775
$ws.Range("M77").Value = -4264.7775

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 85747.46000000001
$ws.Range("I116").Value = 11922.889
$ws.Range("K116").Value = 11922.889
$ws.Range("M116").Value = -9628.888999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5434.697
$ws.Range("I132").Value = 4189.1763
$ws.Range("K132").Value = 12567.5289
$ws.Range("M132").Value = -10037.5289

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 85747.46000000001
$ws.Range("I3").Value = 11922.889
$ws.Range("K3").Value = 11922.889
$ws.Range("M3").Value = -11808.889

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3278.5715
$ws.Range("I20").Value = 2137.5
$ws.Range("J20").Value = 4800
$ws.Range("K20").Value = 2137.5
$ws.Range("L20").Value = 4800
$ws.Range("M20").Value = -1890.5
$ws.Range("N20").Value = -5294

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1261.3125
$ws.Range("I105").Value = 1078.7667
$ws.Range("K105").Value = 1078.7667
$ws.Range("M105").Value = 668.2333000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 17497.414
$ws.Range("I134").Value = 2130.647
$ws.Range("K134").Value = 6391.941
$ws.Range("M134").Value = -3856.941

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 168310.17
$ws.Range("I12").Value = 2505
$ws.Range("K12").Value = 2505
$ws.Range("M12").Value = -2335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 69949.75
$ws.Range("J60").Value = 69942.57000000001
$ws.Range("L60").Value = 69942.57000000001
$ws.Range("N60").Value = -70964.57000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 8953
$ws.Range("J62").Value = 8906
$ws.Range("L62").Value = 8906
$ws.Range("N62").Value = -10154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 8953
$ws.Range("J65").Value = 8906
$ws.Range("L65").Value = 44530
$ws.Range("N65").Value = -50770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8096.75
$ws.Range("I99").Value = 10337
$ws.Range("J99").Value = 7350
$ws.Range("K99").Value = 10337
$ws.Range("L99").Value = 7350
$ws.Range("M99").Value = -8839
$ws.Range("N99").Value = -10346

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 8096.75
$ws.Range("I126").Value = 10337
$ws.Range("J126").Value = 7350
$ws.Range("K126").Value = 31011
$ws.Range("L126").Value = 22050
$ws.Range("M126").Value = -28541
$ws.Range("N126").Value = -26990

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 122855.22
$ws.Range("I134").Value = 1541.3226
$ws.Range("J134").Value = 481020.06
$ws.Range("K134").Value = 4623.9678
$ws.Range("L134").Value = 1443060.18
$ws.Range("M134").Value = -2088.9678
$ws.Range("N134").Value = -1448130.18

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 28.545454
$ws.Range("I38").Value = 29.5
$ws.Range("J38").Value = 19
$ws.Range("K38").Value = 88.5
$ws.Range("L38").Value = 57
$ws.Range("M38").Value = 258.5
$ws.Range("N38").Value = -751

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 50816.617
$ws.Range("J107").Value = 88090.664
$ws.Range("L107").Value = 264271.992
$ws.Range("N107").Value = -268111.992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 905
$ws.Range("J115").Value = 905
$ws.Range("L115").Value = 2715
$ws.Range("N115").Value = -5065

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6439.56
$ws.Range("I70").Value = 6119.533
$ws.Range("J70").Value = 6919.6
$ws.Range("K70").Value = 6119.533
$ws.Range("L70").Value = 6919.6
$ws.Range("M70").Value = -5849.533
$ws.Range("N70").Value = -7459.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6439.56
$ws.Range("I73").Value = 6119.533
$ws.Range("J73").Value = 6919.6
$ws.Range("K73").Value = 6119.533
$ws.Range("L73").Value = 6919.6
$ws.Range("M73").Value = -5183.533
$ws.Range("N73").Value = -8791.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 300947.88
$ws.Range("I132").Value = 419883.88
$ws.Range("K132").Value = 1259651.64
$ws.Range("M132").Value = -1257121.64

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5100
$ws.Range("I46").Value = 4400
$ws.Range("J46").Value = 6500
$ws.Range("K46").Value = 4400
$ws.Range("L46").Value = 6500
$ws.Range("M46").Value = -4212
$ws.Range("N46").Value = -6876

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3966
$ws.Range("I61").Value = 2115.7273
$ws.Range("K61").Value = 2115.7273
$ws.Range("M61").Value = -1913.7273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3359.4
$ws.Range("J82").Value = 1698.75
$ws.Range("L82").Value = 1698.75
$ws.Range("N82").Value = -2420.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3359.4
$ws.Range("J85").Value = 1698.75
$ws.Range("L85").Value = 1698.75
$ws.Range("N85").Value = -4194.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3966
$ws.Range("I113").Value = 2115.7273
$ws.Range("K113").Value = 2115.7273
$ws.Range("M113").Value = 54.27269999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4410.207
$ws.Range("I132").Value = 3866.3809
$ws.Range("J132").Value = 5837.75
$ws.Range("K132").Value = 11599.1427
$ws.Range("L132").Value = 17513.25
$ws.Range("M132").Value = -9069.1427
$ws.Range("N132").Value = -22573.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3561.926
$ws.Range("I136").Value = 3205.6667
$ws.Range("J136").Value = 3846.9333
$ws.Range("K136").Value = 9617.000100000001
$ws.Range("L136").Value = 11540.7999
$ws.Range("M136").Value = -7067.000100000001
$ws.Range("N136").Value = -16640.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 45200
$ws.Range("J93").Value = 45200
$ws.Range("L93").Value = 45200
$ws.Range("N93").Value = -50192

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1959.7333
$ws.Range("I126").Value = 1867.7778
$ws.Range("K126").Value = 5603.3334
$ws.Range("M126").Value = -3133.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 21994.527
$ws.Range("I132").Value = 3179.6428
$ws.Range("K132").Value = 9538.928400000001
$ws.Range("M132").Value = -7008.928400000001
